$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename header cells: "_old" suffix -> "_FV2404", "_new" suffix -> "_FV2410"
# ---------------------------------------------------------------------------
$headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $headers[$i]
}

# ---------------------------------------------------------------------------
# 2. Turn the data range into a proper Excel Table (ListObject) named Table1.
#    The header row already carries its own explicit formatting (bold, grey
#    fill, border). Stash + restore that formatting around the Add() call so
#    the table doesn't record it as a header-row style override (dxf).
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$stashRange = $ws.Range("A100:U100")
$headerRange.Copy()
$stashRange.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$headerRange.ClearFormats()

$range = $ws.Range("A1:U80")
$table = $ws.ListObjects.Add(1, $range, $null, 1)
$table.Name = "Table1"
$table.TableStyle = ""

$stashRange.Copy()
$headerRange.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Rows(100).Delete()

# ---------------------------------------------------------------------------
# 3. Freeze the header row (split pane at row 2).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
